# Switch example dataset from DRC (geomatch) to North America (hmatch)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Wipe the old table (values + formatting) and drop the now-unused 5th
# (E) column so the sheet goes back to a 4-column table.
$ws.Range("A1:E12").Clear()
$ws.Columns.Item(5).Delete()

# New header row: id, adm0, adm1, adm2
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "adm0"
$ws.Range("C1").Value = "adm1"
$ws.Range("D1").Value = "adm2"

# New North-America-themed example rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "united states"
$ws.Range("C2").Value = "new york"
$ws.Range("D2").Value = "suffolk"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Canada"
$ws.Range("C3").Value = "Ontario"

$ws.Range("A4").Value = 3
$ws.Range("D4").Value = "philadelphia"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "United States"
$ws.Range("D5").Value = "York"

$ws.Range("A6").Value = 5
$ws.Range("C6").Value = "NewYork"
$ws.Range("D6").Value = "Jefferson"

$ws.Range("A7").Value = 6
$ws.Range("C7").Value = "pensylvania"
$ws.Range("D7").Value = "philidelphia"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "united_states"
$ws.Range("D8").Value = "king"

$ws.Range("A9").Value = 8
$ws.Range("D9").Value = "NJ_Bergen"

$ws.Range("A10").Value = 9
$ws.Range("D10").Value = "jeffersen"

$ws.Range("A11").Value = 10
$ws.Range("D11").Value = "york"

# Re-apply header styling (centered), matching the original header format
$ws.Range("A1:D1").HorizontalAlignment = -4108

# Autofit-style explicit column widths (B/C/D), chosen so the saved
# width lines up with the widths Excel's own best-fit produced.
$ws.Columns.Item(2).ColumnWidth = 11.498697916666666
$ws.Columns.Item(3).ColumnWidth = 9.830729166666666
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666

# Move the active selection to the last filled cell, like the source edit
$ws.Range("D11").Select() | Out-Null
